$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row for fiscal year 2002 (福島県食肉衛生検査所) is being removed from the
# table; all following rows shift up by one and the used-range dimension
# shrinks from A1:C23 to A1:C22.
$ws.Rows.Item(15).Delete()
